# Update Pomelo price records (D, M, N, O, P, Q, S columns) to match the
# new source data ordering for rows 2-9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; D = 44176; M = 250; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 },
    @{ Row = 3; D = 44491; M = 180; N = 9000;  O = 9000;  P = 9000;  Q = "`$/caja 14 kilos empedrada"; S = 643 },
    @{ Row = 4; D = 44208; M = 210; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 },
    @{ Row = 5; D = 44400; M = 100; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos";           S = 714 },
    @{ Row = 6; D = 44351; M = 300; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 },
    @{ Row = 7; D = 44397; M = 60;  N = 11000; O = 11000; P = 11000; Q = "`$/caja 14 kilos";           S = 786 },
    @{ Row = 8; D = 44309; M = 300; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 },
    @{ Row = 9; D = 44162; M = 120; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("D$row").Value2 = $r.D
    $ws.Range("M$row").Value2 = $r.M
    $ws.Range("N$row").Value2 = $r.N
    $ws.Range("O$row").Value2 = $r.O
    $ws.Range("P$row").Value2 = $r.P
    $ws.Range("Q$row").Value2 = $r.Q
    $ws.Range("S$row").Value2 = $r.S
}
